# Swap the order of the two comma-separated names in the "Recorded By"
# column (G) wherever the value is exactly "System, dnasr281@gmail.com"
# or "admin@admin.com, dnasr281@gmail.com", turning them into
# "dnasr281@gmail.com, System" / "dnasr281@gmail.com, admin@admin.com"
# respectively. All other values are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $val = $cell.Value2

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "admin@admin.com, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, admin@admin.com"
    }
}
